# Swap the contents of column D ("Stn") and column H ("ABS") for rows 1-6,
# so that the ABS column now appears where Stn used to be and vice versa
# (commit: "EXCEL ABS DESPUES STN" -> ABS after/instead-of Stn's old slot).
#
# Row 1 holds plain text headers ("Stn" / "ABS") so a direct value swap is
# enough. Rows 2-6 hold numeric-looking strings (e.g. "338.60") in the cell
# that is moving into column D; Excel would normally reinterpret such text
# as a number (losing the formatting / trailing zero), so those values are
# written with a leading apostrophe to force text storage and then have
# their formatting cleared so no stray number-format style sticks to the
# cell (matching the original "no style" data cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Looks-Numeric($text) {
    return $text -match '^[0-9]+(\.[0-9]+)?$'
}

for ($row = 1; $row -le 6; $row++) {
    $dCell = $ws.Cells.Item($row, 4)  # column D ("Stn" originally)
    $hCell = $ws.Cells.Item($row, 8)  # column H ("ABS" originally)

    $dValue = $dCell.Value2
    $hValue = $hCell.Value2

    # Write H's old value into D.
    if (Looks-Numeric($hValue)) {
        $dCell.Value = "'" + $hValue
        $dCell.ClearFormats()
    } else {
        $dCell.Value = $hValue
    }

    # Write D's old value into H.
    if (Looks-Numeric($dValue)) {
        $hCell.Value = "'" + $dValue
        $hCell.ClearFormats()
    } else {
        $hCell.Value = $dValue
    }
}
